$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the image filename typo (lowercase "image_1353.jpg" -> "IMG_1353.jpg")
# Both E4 and E6 reference the same shared string.
$ws.Range("E4").Value = "IMG_1353.jpg"
$ws.Range("E6").Value = "IMG_1353.jpg"

# Update the active view/selection to cell E6
$ws.Range("E6").Select()
